# Apply the "Updated symbol list" commit:
#  - refresh a handful of Price (column D) quotes in place
#  - insert "One" as the new row 10, pushing WazirX..CoinExToken down one
#    row each (rows 10-18), refreshing their prices and the row-numbered
#    Volume(1h) labels (column E) to match their new row position
#
# All D-column cells (and the re-derived E-column labels) are stored as
# plain text in the source workbook (Price "244.43" etc. is inlineStr,
# not a number), so every write below goes through a leading single-quote
# to force Excel to keep it as text instead of re-typing it as a number,
# and the cell's original Style is restored afterwards so no incidental
# NumberFormat/quote-prefix formatting leaks into the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $savedStyle = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $savedStyle
}

# --- Straight price refreshes (column D only) ---------------------------
Set-TextValue "D2"  "244.51"
Set-TextValue "D3"  "21.95"
Set-TextValue "D4"  "5.448"
Set-TextValue "D5"  "0.05762"
Set-TextValue "D6"  "3.421"
Set-TextValue "D7"  "6.319"
Set-TextValue "D8"  "0.8169"
Set-TextValue "D9"  "1.027"

# --- Rows 10-18: "One" moves to row 10, everything else shifts down -----
Set-TextValue "B10" "One"
Set-TextValue "C10" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.0005842"
Set-TextValue "E10" "9OneONE"

Set-TextValue "B11" "WazirX"
Set-TextValue "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1425"
Set-TextValue "E11" "10WazirXWRX"

Set-TextValue "B12" "MandalaExchangeToken"
Set-TextValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07295"
Set-TextValue "E12" "11MandalaExchangeTokenMDX"

Set-TextValue "B13" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C13" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03130"
Set-TextValue "E13" "12LiechtensteinCryptoassetsExchangeLCX"

Set-TextValue "B14" "BitrueCoin"
Set-TextValue "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03117"
Set-TextValue "E14" "13BitrueCoinBTR"

Set-TextValue "B15" "MCDex"
Set-TextValue "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "4.140"
Set-TextValue "E15" "14MCDexMCB"

Set-TextValue "B16" "BitMartToken"
Set-TextValue "C16" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D16" "0.09371"
Set-TextValue "E16" "15BitMartTokenBMX"

Set-TextValue "B17" "BitForexToken"
Set-TextValue "C17" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001594"
Set-TextValue "E17" "16BitForexTokenBF"

Set-TextValue "B18" "CoinExToken"
Set-TextValue "C18" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04796"
Set-TextValue "E18" "17CoinExTokenCET"

# --- Further straight price refreshes (column D only) -------------------
Set-TextValue "D19" "0.006219"
Set-TextValue "D20" "0.004127"
Set-TextValue "D21" "0.0009907"
Set-TextValue "D23" "3.749"
Set-TextValue "D24" "2.177"
Set-TextValue "D27" "0.0003993"
Set-TextValue "D40" "0.03858"
Set-TextValue "D41" "0.006686"
Set-TextValue "D42" "0.1068"
Set-TextValue "D44" "0.006575"
Set-TextValue "D45" "0.00005612"
Set-TextValue "D47" "0.3895"

Write-Output "Applied symbol list update"
